$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new participant row (row 12) with data uploaded from the frontend file.
$ws.Range("A12").Value = "Lou"
$ws.Range("B12").Value = "Dok"
$ws.Range("C12").Value = "lou@gmail.com"
$ws.Range("D12").Value = "Brunstatt"

# Match the text number-format used by the firstname/lastname columns above.
$ws.Range("A12:B12").NumberFormat = "@"

# Move the active selection the way it ends up after typing the new row (cursor drops to D13).
[void]$ws.Range("D13").Select()
